# Refresh the cryptos list (price + 1h volume change) — GitHub Actions run.
# Source data are plain text cells (inlineStr) in the sheet, so any value
# that merely LOOKS numeric ("297.05", "1.00", ...) must be re-entered with
# a leading apostrophe to keep Excel from auto-converting it to a Number,
# exactly like the original scraped values were stored as text.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.900.85"
$ws.Range("E2").Value = "  -7.59%  "
$ws.Range("D3").Value = "2.520.33"
$ws.Range("E3").Value = "  -3.76%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Value = "'297.05"
$ws.Range("D6").Value = "'93.31"
$ws.Range("E6").Value = "  -7.13%  "
$ws.Range("D7").Value = "'0.570"
$ws.Range("E7").Value = "  -5.43%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.546"
$ws.Range("E9").Value = "  -6.12%  "
$ws.Range("D10").Value = "'35.98"
$ws.Range("E10").Value = "  -8.92%  "
$ws.Range("D11").Value = "'0.0801"
$ws.Range("E11").Value = "  -5.23%  "
$ws.Range("D12").Value = "'7.55"
$ws.Range("E12").Value = "  -7.31%  "
$ws.Range("D13").Value = "'0.107"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "2.906.36"
$ws.Range("E14").Value = "  -3.64%  "
$ws.Range("D15").Value = "2.522.38"
$ws.Range("E15").Value = "  -3.80%  "
$ws.Range("D16").Value = "'0.864"
$ws.Range("E16").Value = "  -6.81%  "
$ws.Range("D17").Value = "'14.02"
$ws.Range("E17").Value = "  -6.67%  "
$ws.Range("D18").Value = "42.898.34"
$ws.Range("E18").Value = "  -7.94%  "
$ws.Range("D19").Value = "0.0₃0963"
$ws.Range("E19").Value = "  -5.01%  "
$ws.Range("D20").Value = "'6.52"
$ws.Range("D21").Value = "'12.28"
$ws.Range("E21").Value = "  -5.48%  "
$ws.Range("D22").Value = "'72.31"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "'258.53"
$ws.Range("E23").Value = "  -6.16%  "
$ws.Range("D24").Value = "'2.90"
$ws.Range("E24").Value = "  -4.81%  "
$ws.Range("D25").Value = "'2.14"
$ws.Range("E25").Value = "  -2.82%  "
$ws.Range("D26").Value = "'28.98"
$ws.Range("E26").Value = "  -1.94%  "
$ws.Range("E27").Value = "  +0.15%  "
$ws.Range("D28").Value = "'9.93"
$ws.Range("E28").Value = "  -6.83%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").Value = "'36.95"
$ws.Range("E29").Value = "  -4.56%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").Value = "'2.13"
$ws.Range("E30").Value = "  -4.13%  "
$ws.Range("D31").Value = "'5.97"
$ws.Range("E31").Value = "  -7.50%  "
$ws.Range("D32").Value = "'3.48"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("E33").Value = "  -2.63%  "
$ws.Range("D34").Value = "'151.22"
$ws.Range("E34").Value = "  -0.75%  "
$ws.Range("D35").Value = "'2.77"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").Value = "'0.0797"
$ws.Range("E36").Value = "  -5.18%  "
$ws.Range("E37").Value = "  -6.59%  "
$ws.Range("E38").Value = "  -3.99%  "
$ws.Range("D39").Value = "'23.83"
$ws.Range("E39").Value = "  +1.54%  "
$ws.Range("D40").Value = "'16.30"
$ws.Range("E40").Value = "  +2.04%  "
$ws.Range("D41").Value = "'3.45"
$ws.Range("E41").Value = "  -5.31%  "
$ws.Range("E42").Value = "  -7.22%  "
$ws.Range("E43").Value = "  -5.93%  "
$ws.Range("D44").Value = "2.016.08"
$ws.Range("E44").Value = "  -5.49%  "
$ws.Range("D45").Value = "'0.997"
$ws.Range("E45").Value = "  -0.16%  "
$ws.Range("D46").Value = "'85.37"
$ws.Range("E46").Value = "  -9.48%  "
$ws.Range("E47").Value = "  +3.31%  "
$ws.Range("E48").Value = "  -7.19%  "
$ws.Range("D49").Value = "2.763.71"
$ws.Range("E49").Value = "  -3.83%  "
$ws.Range("D50").Value = "'102.64"
$ws.Range("E50").Value = "  -6.74%  "
$ws.Range("E51").Value = "  -8.30%  "
